# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Updates the "Valor Mora" total and "Cant. Periodos" count
# - Inserts a new data row (row 25) into the worker table, preserving the
#   special bottom-border formatting on the new last row
# - Re-populates all data rows (16-25) with the refreshed period/value data:
#     EDGAR JOSE ROA AMADOR (CC 1047386377): periods 2412,2501-2508
#     MAVINETH CECILIA LOPEZ BERMUDEZ (CC 45767046): period 2504 (mora row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row before row 25, shifting rows 25+ down by one ---
$ws.Rows("25:25").Insert()

# --- 2. Fix up formatting: row 24 currently still carries the "last row"
#        (bottom-border) styling that belongs on the new final row 25 now ---
$ws.Range("B24:J24").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B23:J23").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Update header summary cells ---
$ws.Range("E11").Value = 604266            # VALOR MORA total
$ws.Range("F13").Value = 9                 # Cant. Periodos

# --- 4. Re-write the detail table rows (worker/period/value) ---
# Row 16: EDGAR - 2412 (partial period)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047386377"
$ws.Range("D16").Value = "EDGAR JOSE ROA AMADOR"
$ws.Range("E16").Value = "2412"
$ws.Range("F16").Value = 24266
$ws.Range("G16").Value = 1300000

# Row 17: EDGAR - 2501
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047386377"
$ws.Range("D17").Value = "EDGAR JOSE ROA AMADOR"
$ws.Range("E17").Value = "2501"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

# Row 18: EDGAR - 2502
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047386377"
$ws.Range("D18").Value = "EDGAR JOSE ROA AMADOR"
$ws.Range("E18").Value = "2502"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

# Row 19: EDGAR - 2503
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047386377"
$ws.Range("D19").Value = "EDGAR JOSE ROA AMADOR"
$ws.Range("E19").Value = "2503"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

# Row 20: EDGAR - 2504
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047386377"
$ws.Range("D20").Value = "EDGAR JOSE ROA AMADOR"
$ws.Range("E20").Value = "2504"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000

# Row 21: MAVINETH - 2504
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "45767046"
$ws.Range("D21").Value = "MAVINETH CECILIA LOPEZ BERMUDEZ"
$ws.Range("E21").Value = "2504"
$ws.Range("F21").Value = 164000
$ws.Range("G21").Value = 4545527

# Row 22: EDGAR - 2505
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047386377"
$ws.Range("D22").Value = "EDGAR JOSE ROA AMADOR"
$ws.Range("E22").Value = "2505"
$ws.Range("F22").Value = 52000
$ws.Range("G22").Value = 1300000

# Row 23: EDGAR - 2506
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1047386377"
$ws.Range("D23").Value = "EDGAR JOSE ROA AMADOR"
$ws.Range("E23").Value = "2506"
$ws.Range("F23").Value = 52000
$ws.Range("G23").Value = 1300000

# Row 24: EDGAR - 2507
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1047386377"
$ws.Range("D24").Value = "EDGAR JOSE ROA AMADOR"
$ws.Range("E24").Value = "2507"
$ws.Range("F24").Value = 52000
$ws.Range("G24").Value = 1300000

# Row 25 (new): EDGAR - 2508
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1047386377"
$ws.Range("D25").Value = "EDGAR JOSE ROA AMADOR"
$ws.Range("E25").Value = "2508"
$ws.Range("F25").Value = 52000
$ws.Range("G25").Value = 1300000
